$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 192: add MOVE value, fix VIX value
$ws.Range("B192").Value = 60.61
$ws.Range("C192").Value = 22.96

# New data rows 193-198 (Serie date text, MOVE, VIX)
$data = @(
    @("05-10-2021", 62.83, 21.3),
    @("06-10-2021", 60.65, 21),
    @("07-10-2021", 60.91, 19.54),
    @("08-10-2021", 59.65, 18.77),
    @("11-10-2021", $null, 20),
    @("12-10-2021", $null, 19.56)
)

$row = 193
foreach ($entry in $data) {
    $cellA = $ws.Cells.Item($row, 1)
    # Use a formula that evaluates to the literal text, then paste the
    # resulting value back in place so the date-like string is stored as
    # plain text (shared string) rather than being auto-converted to a
    # serial date value / acquiring a new number-format style.
    $cellA.Formula = '="' + $entry[0] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    if ($null -ne $entry[1]) {
        $ws.Cells.Item($row, 2).Value = $entry[1]
    }
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
